$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("K") values were regenerated to reflect strikeouts (K) instead
# of the old "Strike#" metric. Update the four data rows accordingly.
$ws.Range("G2").Value = 4
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 4
$ws.Range("G5").Value = 4
